$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data: the Ip value in row 2 changed ---
$ws.Range("A2").Value = "10.16.5.64"

# --- Add the two new columns (BC, BD) for the new "Decisión" / "Transacción 06-2001" data ---

# Header row (row 1): BC1 = "Decisión", BD1 = "Transacción"
$ws.Range("BC1").Value = "Decisión"
$ws.Range("BD1").Value = "Transacción"

# Copy header formatting (style used by the rest of row 1) onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("BC1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("BD1").PasteSpecial(-4122)

# Data row (row 2): BC2 = "APROBAR", BD2 = 2001 (numeric)
$ws.Range("BC2").Value = "APROBAR"
$ws.Range("BD2").Value = 2001

# Copy the data-row formatting onto the new data cells (BC2 matches the text cells in row 2)
$ws.Range("B2").Copy()
$ws.Range("BC2").PasteSpecial(-4122)
$ws.Range("T2").Copy()
$ws.Range("BD2").PasteSpecial(-4122)
